# Updated symbol list on Fri Dec 16 23:05:36 UTC 2022 with GitHub Actions
#
# The "Price" (D) and "Hora" (G) columns store numeric-looking values as
# text in this workbook. Mark the cells we are about to rewrite as Text
# first so Excel does not silently re-interpret the strings as numbers
# (only the cells whose value actually changes are touched, to avoid
# stamping an unrelated number format onto untouched cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCells = @( `
    "D2","D3","D4","D5","D6","D8","D9","D10","D11","D12","D13","D14", `
    "D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D26","D27", `
    "D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51" `
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("G2:G51").NumberFormat = "@"

# --- Row 2 (BNB) ---
$ws.Range("D2").Value = '231.62'
$ws.Range("G2").Value = '23'

# --- Row 3 (OKB) ---
$ws.Range("D3").Value = '22.89'
$ws.Range("G3").Value = '23'

# --- Row 4 (HuobiToken) ---
$ws.Range("D4").Value = '5.523'
$ws.Range("G4").Value = '23'

# --- Row 5 (Cronos) ---
$ws.Range("D5").Value = '0.05614'
$ws.Range("G5").Value = '23'

# --- Row 6 (GateToken) ---
$ws.Range("D6").Value = '3.418'
$ws.Range("G6").Value = '23'

# --- Row 7 (KuCoinToken) ---
$ws.Range("G7").Value = '23'

# --- Row 8 (FTXToken) ---
$ws.Range("D8").Value = '1.245'
$ws.Range("G8").Value = '23'

# --- Row 9 (MXToken) ---
$ws.Range("D9").Value = '0.8009'
$ws.Range("G9").Value = '23'

# --- Row 10 (WazirX) ---
$ws.Range("D10").Value = '0.1424'
$ws.Range("G10").Value = '23'

# --- Row 11 (MandalaExchangeToken) ---
$ws.Range("D11").Value = '0.07475'
$ws.Range("G11").Value = '23'

# --- Row 12 (LiechtensteinCryptoassetsExchange) ---
$ws.Range("D12").Value = '0.03174'
$ws.Range("G12").Value = '23'

# --- Row 13 (BitrueCoin) ---
$ws.Range("D13").Value = '0.02942'
$ws.Range("G13").Value = '23'

# --- Row 14 (BitMartToken) ---
$ws.Range("D14").Value = '0.09248'
$ws.Range("G14").Value = '23'

# --- Row 15 (BitForexToken) ---
$ws.Range("D15").Value = '0.001679'
$ws.Range("G15").Value = '23'

# --- Row 16 (MCDex) ---
$ws.Range("D16").Value = '3.268'
$ws.Range("G16").Value = '23'

# --- Row 17 (CoinExToken) ---
$ws.Range("D17").Value = '0.04731'
$ws.Range("G17").Value = '23'

# --- Row 18: was "One", now "TigerCash" (rows 18-24 shift up one rank) ---
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '0.006236'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("G18").Value = '23'

# --- Row 19: was "TigerCash", now "HotbitToken" ---
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").Value = '0.005311'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("G19").Value = '23'

# --- Row 20: was "HotbitToken", now "BitKan" ---
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").Value = '0.001067'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("G20").Value = '23'

# --- Row 21: was "BitKan", now "NitroEx" ---
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Value = '0.0001506'
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("G21").Value = '23'

# --- Row 22: was "NitroEx", now "LEO" ---
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").Value = '3.686'
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("G22").Value = '23'

# --- Row 23: was "LEO", now "BTSEToken" ---
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").Value = '2.192'
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("G23").Value = '23'

# --- Row 24: was "BTSEToken", now "One" ---
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").Value = '0.0005988'
$ws.Range("E24").Value = '23OneONE'
$ws.Range("G24").Value = '23'

# --- Row 25 (BitpandaEcosystemToken) ---
$ws.Range("G25").Value = '23'

# --- Row 26 (ProBitToken) ---
$ws.Range("D26").Value = '0.1281'
$ws.Range("G26").Value = '23'

# --- Row 27 (UpBots) ---
$ws.Range("D27").Value = '0.0006663'
$ws.Range("G27").Value = '23'

# --- Rows 28-39: price unchanged ("--"), only Hora changes ---
$ws.Range("G28").Value = '23'
$ws.Range("G29").Value = '23'
$ws.Range("G30").Value = '23'
$ws.Range("G31").Value = '23'
$ws.Range("G32").Value = '23'
$ws.Range("G33").Value = '23'
$ws.Range("G34").Value = '23'
$ws.Range("G35").Value = '23'
$ws.Range("G36").Value = '23'
$ws.Range("G37").Value = '23'
$ws.Range("G38").Value = '23'
$ws.Range("G39").Value = '23'

# --- Row 40 (IDEX) ---
$ws.Range("D40").Value = '0.04121'
$ws.Range("G40").Value = '23'

# --- Row 41 (KickToken) ---
$ws.Range("D41").Value = '0.007128'
$ws.Range("G41").Value = '23'

# --- Row 42 (CEJI) ---
$ws.Range("D42").Value = '0.003459'
$ws.Range("G42").Value = '23'

# --- Row 43 (BKEXToken) ---
$ws.Range("D43").Value = '0.1042'
$ws.Range("G43").Value = '23'

# --- Row 44 (LocalTraders) ---
$ws.Range("D44").Value = '0.009257'
$ws.Range("G44").Value = '23'

# --- Row 45 (ACDXExchange) ---
$ws.Range("E45").Value = '44ACDXExchangeACXTWorstin24h'
$ws.Range("G45").Value = '23'

# --- Row 46 (CoinLion) ---
$ws.Range("D46").Value = '0.00005592'
$ws.Range("G46").Value = '23'

# --- Row 47 (Kangarootoken) ---
$ws.Range("D47").Value = '0.00000000753'
$ws.Range("G47").Value = '23'

# --- Row 48 (CoinbaseStockToken) ---
$ws.Range("D48").Value = '0.7888'
$ws.Range("G48").Value = '23'

# --- Row 49 (BOLO) ---
$ws.Range("D49").Value = '0.09599'
$ws.Range("G49").Value = '23'

# --- Row 50 (CryptobidCoin) ---
$ws.Range("D50").Value = '0.00002109'
$ws.Range("G50").Value = '23'

# --- Row 51 (SpecialPowerGold) ---
$ws.Range("D51").Value = '0.01014'
$ws.Range("G51").Value = '23'
